$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")
$ws.Columns("C:C").Insert()
$ws.Range("C1").Value = "appearance"
$ws.Range("C5").Value = "grid"
$ws.Range("B9").Value = "begin screen"
$ws.Range("B10").Value = "select_one yes_no"
$ws.Range("C10").Value = "inline"
$ws.Range("E10").Value = "i1"
$ws.Range("F10").Value = "Choose one:"
$ws.Range("B11").Value = "select_one yes_no"
$ws.Range("C11").Value = "inline"
$ws.Range("E11").Value = "i2"
$ws.Range("F11").Value = "Choose one:"
$ws.Range("B12").Value = "select_one yes_no"
$ws.Range("C12").Value = "inline"
$ws.Range("E12").Value = "i3"
$ws.Range("F12").Value = "Choose one:"
$ws.Range("B13").Value = "end screen"

$ws2 = $wb.Worksheets.Item("choices")
$ws2.Range("D1").Value = "label"
$ws2.Range("A12").Value = "yes_no"
$ws2.Range("B12").Value = "yes"
$ws2.Range("D12").Value = "Yes"
$ws2.Range("A13").Value = "yes_no"
$ws2.Range("B13").Value = "no"
$ws2.Range("D13").Value = "No"

$ws3 = $wb.Worksheets.Item("queries")
$ws3.Range("D1").Clear()
$ws3.Range("E1").Clear()
$ws3.Range("B5").Value = """content://com.opendatakit.tables.ContentProvider/database_id/table_id/row_id"""
